$d = $word.ActiveDocument

# The "nibble1" bit-field row (Bits [11:8]) has a Description cell that
# currently reads "write something useful for nibble1". Clear it out so
# the cell is left with a single, completely empty paragraph - same as
# the sibling "nibble1" row elsewhere in the document whose description
# was already left blank.
$t = $d.Tables.Item(13)
$cell = $t.Cell(5, 4)

$cell.Range.Find.Execute("write something useful for nibble1", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$para = $cell.Range.Paragraphs.Item(1)
$para.Style = "Normal"
